$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 533.5
$ws.Range("I4").Value = 175.5
$ws.Range("J4").Value = 1249.5
$ws.Range("K4").Value = 175.5
$ws.Range("L4").Value = 1249.5
$ws.Range("M4").Value = -61.5
$ws.Range("N4").Value = -1477.5
$ws.Range("H10").Value = 27500
$ws.Range("I10").Value = 20000
$ws.Range("K10").Value = 20000
$ws.Range("M10").Value = -19707
$ws.Range("H15").Value = 1157
$ws.Range("I15").Value = 1157
$ws.Range("K15").Value = 3471
$ws.Range("M15").Value = -3302
$ws.Range("H17").Value = 2100
$ws.Range("J17").Value = 2320
$ws.Range("L17").Value = 6960
$ws.Range("N17").Value = -7296
$ws.Range("H28").Value = 315.2143
$ws.Range("I28").Value = 354.25
$ws.Range("K28").Value = 354.25
$ws.Range("M28").Value = 130.75
$ws.Range("H32").Value = 2057.9092
$ws.Range("I32").Value = 2459.8
$ws.Range("J32").Value = 1723
$ws.Range("K32").Value = 2459.8
$ws.Range("L32").Value = 1723
$ws.Range("M32").Value = -2133.8
$ws.Range("N32").Value = -2375
$ws.Range("H40").Value = 1539.8572
$ws.Range("I40").Value = 1435.8
$ws.Range("K40").Value = 1435.8
$ws.Range("M40").Value = -1260.8
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2138
$ws.Range("H69").Value = 6680.25
$ws.Range("I69").Value = 4013
$ws.Range("J69").Value = 6922.727
$ws.Range("K69").Value = 12039
$ws.Range("L69").Value = 20768.181
$ws.Range("M69").Value = -11165
$ws.Range("N69").Value = -22516.181
$ws.Range("H70").Value = 4010.125
$ws.Range("J70").Value = 3940.2856
$ws.Range("L70").Value = 11820.8568
$ws.Range("N70").Value = -12360.8568
$ws.Range("H72").Value = 6680.25
$ws.Range("I72").Value = 4013
$ws.Range("J72").Value = 6922.727
$ws.Range("K72").Value = 36117
$ws.Range("L72").Value = 62304.543
$ws.Range("M72").Value = -31749
$ws.Range("N72").Value = -71040.54300000001
$ws.Range("H73").Value = 4010.125
$ws.Range("J73").Value = 3940.2856
$ws.Range("L73").Value = 11820.8568
$ws.Range("N73").Value = -13692.8568
$ws.Range("H76").Value = 3155.8
$ws.Range("J76").Value = 3914.2307
$ws.Range("L76").Value = 3914.2307
$ws.Range("N76").Value = -4544.2307
$ws.Range("H79").Value = 3155.8
$ws.Range("J79").Value = 3914.2307
$ws.Range("L79").Value = 3914.2307
$ws.Range("N79").Value = -6098.2307
$ws.Range("H80").Value = 2119.5417
$ws.Range("I80").Value = 2220.4
$ws.Range("J80").Value = 2093
$ws.Range("K80").Value = 6661.200000000001
$ws.Range("L80").Value = 6279
$ws.Range("M80").Value = -5663.200000000001
$ws.Range("N80").Value = -8275
$ws.Range("H83").Value = 2119.5417
$ws.Range("I83").Value = 2220.4
$ws.Range("J83").Value = 2093
$ws.Range("K83").Value = 19983.6
$ws.Range("L83").Value = 18837
$ws.Range("M83").Value = -14991.6
$ws.Range("N83").Value = -28821
$ws.Range("H86").Value = 3180.7856
$ws.Range("I86").Value = 2859
$ws.Range("K86").Value = 2859
$ws.Range("M86").Value = -1736
$ws.Range("H87").Value = 63725
$ws.Range("J87").Value = 63725
$ws.Range("L87").Value = 63725
$ws.Range("N87").Value = -66221
$ws.Range("H88").Value = 1243.8148
$ws.Range("I88").Value = 1324
$ws.Range("J88").Value = 1203.7222
$ws.Range("K88").Value = 1324
$ws.Range("L88").Value = 1203.7222
$ws.Range("M88").Value = -918
$ws.Range("N88").Value = -2015.7222
$ws.Range("H89").Value = 3180.7856
$ws.Range("I89").Value = 2859
$ws.Range("K89").Value = 14295
$ws.Range("M89").Value = -8679
$ws.Range("H90").Value = 63725
$ws.Range("J90").Value = 63725
$ws.Range("L90").Value = 191175
$ws.Range("N90").Value = -203655
$ws.Range("H91").Value = 1243.8148
$ws.Range("I91").Value = 1324
$ws.Range("J91").Value = 1203.7222
$ws.Range("K91").Value = 1324
$ws.Range("L91").Value = 1203.7222
$ws.Range("M91").Value = 80
$ws.Range("N91").Value = -4011.7222
$ws.Range("H106").Value = 2482.1875
$ws.Range("I106").Value = 2793.4614
$ws.Range("K106").Value = 2793.4614
$ws.Range("M106").Value = -2162.4614
$ws.Range("H111").Value = 1363.7778
$ws.Range("I111").Value = 1258.1666
$ws.Range("J111").Value = 1575
$ws.Range("K111").Value = 3774.4998
$ws.Range("L111").Value = 4725
$ws.Range("M111").Value = -707.4998000000001
$ws.Range("N111").Value = -10859
$ws.Range("H112").Value = 2766.4644
$ws.Range("I112").Value = 1925.7
$ws.Range("J112").Value = 3233.5557
$ws.Range("K112").Value = 5777.1
$ws.Range("L112").Value = 9700.667099999999
$ws.Range("M112").Value = -4669.1
$ws.Range("N112").Value = -11916.6671
$ws.Range("H116").Value = 18890.625
$ws.Range("I116").Value = 4639
$ws.Range("K116").Value = 4639
$ws.Range("M116").Value = -1197
$ws.Range("H131").Value = 14886.454
$ws.Range("I131").Value = 1375.1
$ws.Range("K131").Value = 4125.299999999999
$ws.Range("M131").Value = 914.7000000000007
$ws.Range("H138").Value = 3605.7666
$ws.Range("I138").Value = 6043
$ws.Range("J138").Value = 3431.6785
$ws.Range("K138").Value = 18129
$ws.Range("L138").Value = 10295.0355
$ws.Range("M138").Value = -12989
$ws.Range("N138").Value = -20575.0355
$ws.Range("H140").Value = 104206.93
$ws.Range("J140").Value = 104206.93
$ws.Range("L140").Value = 104206.93
$ws.Range("N140").Value = -114566.93
$ws.Range("H141").Value = 5175.2104
$ws.Range("I141").Value = 5296.0557
$ws.Range("K141").Value = 15888.1671
$ws.Range("M141").Value = -10708.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1796.4445
$ws.Range("I5").Value = 159.14285
$ws.Range("J5").Value = 7527
$ws.Range("K5").Value = 159.14285
$ws.Range("L5").Value = 7527
$ws.Range("M5").Value = -47.14285000000001
$ws.Range("N5").Value = -7751
$ws.Range("H32").Value = 150954.45
$ws.Range("I32").Value = 197448.48
$ws.Range("J32").Value = 19221.334
$ws.Range("K32").Value = 197448.48
$ws.Range("L32").Value = 19221.334
$ws.Range("M32").Value = -197161.48
$ws.Range("N32").Value = -19795.334
$ws.Range("H45").Value = 2154.3333
$ws.Range("I45").Value = 2278
$ws.Range("K45").Value = 2278
$ws.Range("M45").Value = -1901
$ws.Range("H74").Value = 1859329.6
$ws.Range("I74").Value = 3971346.8
$ws.Range("K74").Value = 3971346.8
$ws.Range("M74").Value = -3970472.8
$ws.Range("H76").Value = 27950
$ws.Range("J76").Value = 27950
$ws.Range("L76").Value = 27950
$ws.Range("N76").Value = -28626
$ws.Range("H77").Value = 1859329.6
$ws.Range("I77").Value = 3971346.8
$ws.Range("K77").Value = 19856734
$ws.Range("M77").Value = -19852366
$ws.Range("H79").Value = 27950
$ws.Range("J79").Value = 27950
$ws.Range("L79").Value = 27950
$ws.Range("N79").Value = -30290
$ws.Range("H97").Value = 1343.4
$ws.Range("I97").Value = 1649
$ws.Range("J97").Value = 630.3333
$ws.Range("K97").Value = 1649
$ws.Range("L97").Value = 630.3333
$ws.Range("M97").Value = -1153
$ws.Range("N97").Value = -1622.3333
$ws.Range("H122").Value = 1164.0605
$ws.Range("I122").Value = 1042.4138
$ws.Range("J122").Value = 2046
$ws.Range("K122").Value = 3127.2414
$ws.Range("L122").Value = 6138
$ws.Range("M122").Value = -677.2413999999999
$ws.Range("N122").Value = -11038
$ws.Range("H132").Value = 512504.1
$ws.Range("I132").Value = 570467.4399999999
$ws.Range("J132").Value = 2426.6
$ws.Range("K132").Value = 1711402.32
$ws.Range("L132").Value = 7279.799999999999
$ws.Range("M132").Value = -1708872.32
$ws.Range("N132").Value = -12339.8
$ws.Range("H133").Value = 84034.86
$ws.Range("J133").Value = 84034.86
$ws.Range("L133").Value = 84034.86
$ws.Range("N133").Value = -89094.86
$ws.Range("H140").Value = 61942.332
$ws.Range("J140").Value = 61942.332
$ws.Range("L140").Value = 61942.332
$ws.Range("N140").Value = -72302.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1796.4445
$ws.Range("I4").Value = 159.14285
$ws.Range("J4").Value = 7527
$ws.Range("K4").Value = 159.14285
$ws.Range("L4").Value = 7527
$ws.Range("M4").Value = -44.14285000000001
$ws.Range("N4").Value = -7757
$ws.Range("H82").Value = 23950.8
$ws.Range("I82").Value = 10501.333
$ws.Range("J82").Value = 44125
$ws.Range("K82").Value = 10501.333
$ws.Range("L82").Value = 44125
$ws.Range("M82").Value = -10118.333
$ws.Range("N82").Value = -44891
$ws.Range("H85").Value = 23950.8
$ws.Range("I85").Value = 10501.333
$ws.Range("J85").Value = 44125
$ws.Range("K85").Value = 10501.333
$ws.Range("L85").Value = 44125
$ws.Range("M85").Value = -9175.333000000001
$ws.Range("N85").Value = -46777
$ws.Range("H134").Value = 4171339
$ws.Range("I134").Value = 4743.4688
$ws.Range("K134").Value = 14230.4064
$ws.Range("M134").Value = -11695.4064
$ws.Range("H140").Value = 79292.664
$ws.Range("J140").Value = 79292.664
$ws.Range("L140").Value = 79292.664
$ws.Range("N140").Value = -89652.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47657.523
$ws.Range("I7").Value = 62527.875
$ws.Range("J7").Value = 72.40000000000001
$ws.Range("K7").Value = 62527.875
$ws.Range("L7").Value = 72.40000000000001
$ws.Range("M7").Value = -62414.875
$ws.Range("N7").Value = -298.4
$ws.Range("H16").Value = 2866.7646
$ws.Range("I16").Value = 2075.7693
$ws.Range("K16").Value = 2075.7693
$ws.Range("M16").Value = -1788.7693
$ws.Range("H31").Value = 1545515.8
$ws.Range("I31").Value = 2060011.5
$ws.Range("J31").Value = 2028.5555
$ws.Range("K31").Value = 2060011.5
$ws.Range("L31").Value = 2028.5555
$ws.Range("M31").Value = -2059716.5
$ws.Range("N31").Value = -2618.5555
$ws.Range("H34").Value = 1545515.8
$ws.Range("I34").Value = 2060011.5
$ws.Range("J34").Value = 2028.5555
$ws.Range("K34").Value = 2060011.5
$ws.Range("L34").Value = 2028.5555
$ws.Range("M34").Value = -2059809.5
$ws.Range("N34").Value = -2432.5555
$ws.Range("H58").Value = 2529194.2
$ws.Range("I58").Value = 3245.4666
$ws.Range("J58").Value = 4634151.5
$ws.Range("K58").Value = 3245.4666
$ws.Range("L58").Value = 4634151.5
$ws.Range("M58").Value = -3042.4666
$ws.Range("N58").Value = -4634557.5
$ws.Range("H62").Value = 3320.5
$ws.Range("J62").Value = 3316.25
$ws.Range("L62").Value = 3316.25
$ws.Range("N62").Value = -4564.25
$ws.Range("H65").Value = 3320.5
$ws.Range("J65").Value = 3316.25
$ws.Range("L65").Value = 16581.25
$ws.Range("N65").Value = -22821.25
$ws.Range("H92").Value = 85000
$ws.Range("J92").Value = 85000
$ws.Range("L92").Value = 85000
$ws.Range("N92").Value = -89992
$ws.Range("H94").Value = 1760.7646
$ws.Range("J94").Value = 1918.2307
$ws.Range("L94").Value = 1918.2307
$ws.Range("N94").Value = -2820.2307
$ws.Range("H99").Value = 22264.875
$ws.Range("I99").Value = 22204
$ws.Range("J99").Value = 22447.5
$ws.Range("K99").Value = 22204
$ws.Range("L99").Value = 22447.5
$ws.Range("M99").Value = -20706
$ws.Range("N99").Value = -25443.5
$ws.Range("H102").Value = 30097
$ws.Range("I102").Value = 23797
$ws.Range("J102").Value = 38497
$ws.Range("K102").Value = 23797
$ws.Range("L102").Value = 38497
$ws.Range("M102").Value = -21363
$ws.Range("N102").Value = -43365
$ws.Range("H105").Value = 1574.6875
$ws.Range("I105").Value = 1146.4166
$ws.Range("K105").Value = 1146.4166
$ws.Range("M105").Value = 600.5834
$ws.Range("H107").Value = 560.8
$ws.Range("I107").Value = 559.8333
$ws.Range("J107").Value = 562.25
$ws.Range("K107").Value = 559.8333
$ws.Range("L107").Value = 562.25
$ws.Range("M107").Value = 1360.1667
$ws.Range("N107").Value = -4402.25
$ws.Range("H113").Value = 2866.7646
$ws.Range("I113").Value = 2075.7693
$ws.Range("K113").Value = 2075.7693
$ws.Range("M113").Value = 94.23070000000007
$ws.Range("H126").Value = 22264.875
$ws.Range("I126").Value = 22204
$ws.Range("J126").Value = 22447.5
$ws.Range("K126").Value = 66612
$ws.Range("L126").Value = 67342.5
$ws.Range("M126").Value = -64142
$ws.Range("N126").Value = -72282.5
$ws.Range("H129").Value = 41246.668
$ws.Range("I129").Value = 33800
$ws.Range("K129").Value = 33800
$ws.Range("M129").Value = -28800
$ws.Range("H132").Value = 2534.2195
$ws.Range("I132").Value = 2500.1282
$ws.Range("J132").Value = 3199
$ws.Range("K132").Value = 7500.3846
$ws.Range("L132").Value = 9597
$ws.Range("M132").Value = -4970.3846
$ws.Range("N132").Value = -14657
$ws.Range("H136").Value = 2529194.2
$ws.Range("I136").Value = 3245.4666
$ws.Range("J136").Value = 4634151.5
$ws.Range("K136").Value = 9736.399800000001
$ws.Range("L136").Value = 13902454.5
$ws.Range("M136").Value = -7186.399800000001
$ws.Range("N136").Value = -13907554.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4915429.5
$ws.Range("I5").Value = 7143617.5
$ws.Range("J5").Value = 3323866.5
$ws.Range("K5").Value = 21430852.5
$ws.Range("L5").Value = 9971599.5
$ws.Range("M5").Value = -21430740.5
$ws.Range("N5").Value = -9971823.5
$ws.Range("H17").Value = 1497.5
$ws.Range("I17").Value = 295
$ws.Range("J17").Value = 2700
$ws.Range("K17").Value = 885
$ws.Range("L17").Value = 8100
$ws.Range("M17").Value = -716
$ws.Range("N17").Value = -8438
$ws.Range("H64").Value = 9178
$ws.Range("I64").Value = 8750
$ws.Range("J64").Value = 9249.333000000001
$ws.Range("K64").Value = 26250
$ws.Range("L64").Value = 27747.999
$ws.Range("M64").Value = -25980
$ws.Range("N64").Value = -28287.999
$ws.Range("H67").Value = 9178
$ws.Range("I67").Value = 8750
$ws.Range("J67").Value = 9249.333000000001
$ws.Range("K67").Value = 26250
$ws.Range("L67").Value = 27747.999
$ws.Range("M67").Value = -25314
$ws.Range("N67").Value = -29619.999
$ws.Range("H69").Value = 2999.75
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 2999.75
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112
$ws.Range("H74").Value = 14583.333
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 14583.333
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 43749.999
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -45871.999
$ws.Range("H77").Value = 14583.333
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 14583.333
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 131249.997
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -141857.997
$ws.Range("H80").Value = 11400
$ws.Range("J80").Value = 11400
$ws.Range("L80").Value = 34200
$ws.Range("N80").Value = -36072
$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8246
$ws.Range("H82").Value = 13909.091
$ws.Range("J82").Value = 13909.091
$ws.Range("L82").Value = 41727.273
$ws.Range("N82").Value = -42539.273
$ws.Range("H83").Value = 11400
$ws.Range("J83").Value = 11400
$ws.Range("L83").Value = 102600
$ws.Range("N83").Value = -111960
$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 18000
$ws.Range("N84").Value = -29232
$ws.Range("H85").Value = 13909.091
$ws.Range("J85").Value = 13909.091
$ws.Range("L85").Value = 41727.273
$ws.Range("N85").Value = -44535.273
$ws.Range("H132").Value = 1523.4
$ws.Range("I132").Value = 999.5
$ws.Range("J132").Value = 1654.375
$ws.Range("K132").Value = 8995.5
$ws.Range("L132").Value = 14889.375
$ws.Range("M132").Value = -6465.5
$ws.Range("N132").Value = -19949.375
$ws.Range("H135").Value = 4915429.5
$ws.Range("I135").Value = 7143617.5
$ws.Range("J135").Value = 3323866.5
$ws.Range("K135").Value = 64292557.5
$ws.Range("L135").Value = 29914798.5
$ws.Range("M135").Value = -64290022.5
$ws.Range("N135").Value = -29919868.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 15000
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15554
$ws.Range("H62").Value = 62000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 62000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 62000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -63372
$ws.Range("H65").Value = 62000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 62000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 186000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -192864
$ws.Range("H70").Value = 33284.42
$ws.Range("I70").Value = 8173.7827
$ws.Range("K70").Value = 8173.7827
$ws.Range("M70").Value = -7903.7827
$ws.Range("H73").Value = 33284.42
$ws.Range("I73").Value = 8173.7827
$ws.Range("K73").Value = 8173.7827
$ws.Range("M73").Value = -7237.7827
$ws.Range("H80").Value = 4671.826
$ws.Range("I80").Value = 3176.353
$ws.Range("K80").Value = 3176.353
$ws.Range("M80").Value = -2178.353
$ws.Range("H83").Value = 4671.826
$ws.Range("I83").Value = 3176.353
$ws.Range("K83").Value = 15881.765
$ws.Range("M83").Value = -10889.765
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -22744
$ws.Range("H107").Value = 916.0606
$ws.Range("I107").Value = 705.2083
$ws.Range("K107").Value = 705.2083
$ws.Range("M107").Value = 1214.7917
$ws.Range("H113").Value = 2731.1785
$ws.Range("I113").Value = 2687.1333
$ws.Range("J113").Value = 2782
$ws.Range("K113").Value = 2687.1333
$ws.Range("L113").Value = 2782
$ws.Range("M113").Value = -517.1333
$ws.Range("N113").Value = -7122
$ws.Range("H122").Value = 42940.73
$ws.Range("I122").Value = 53733.55
$ws.Range("K122").Value = 161200.65
$ws.Range("M122").Value = -158750.65
$ws.Range("H124").Value = 84999.5
$ws.Range("J124").Value = 84999.5
$ws.Range("L124").Value = 84999.5
$ws.Range("N124").Value = -94819.5
$ws.Range("H132").Value = 2681.1445
$ws.Range("I132").Value = 2184.0422
$ws.Range("K132").Value = 6552.1266
$ws.Range("M132").Value = -4022.1266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1698.5834
$ws.Range("I16").Value = 427
$ws.Range("K16").Value = 427
$ws.Range("M16").Value = -257
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H40").Value = 4042.0667
$ws.Range("I40").Value = 4188.2
$ws.Range("K40").Value = 4188.2
$ws.Range("M40").Value = -4052.2
$ws.Range("H46").Value = 3968.375
$ws.Range("I46").Value = 493.6
$ws.Range("J46").Value = 5547.8184
$ws.Range("K46").Value = 493.6
$ws.Range("L46").Value = 5547.8184
$ws.Range("M46").Value = -305.6
$ws.Range("N46").Value = -5923.8184
$ws.Range("H61").Value = 7730.5757
$ws.Range("I61").Value = 6675.077
$ws.Range("J61").Value = 11651
$ws.Range("K61").Value = 6675.077
$ws.Range("L61").Value = 11651
$ws.Range("M61").Value = -6473.077
$ws.Range("N61").Value = -12055
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H68").Value = 1892.3077
$ws.Range("I68").Value = 1892.3077
$ws.Range("K68").Value = 1892.3077
$ws.Range("M68").Value = -1143.3077
$ws.Range("H71").Value = 1892.3077
$ws.Range("I71").Value = 1892.3077
$ws.Range("K71").Value = 9461.538500000001
$ws.Range("M71").Value = -5717.538500000001
$ws.Range("H82").Value = 1905.1818
$ws.Range("I82").Value = 2751.4
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 2751.4
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -2390.4
$ws.Range("N82").Value = -1922
$ws.Range("H85").Value = 1905.1818
$ws.Range("I85").Value = 2751.4
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 2751.4
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -1503.4
$ws.Range("N85").Value = -3696
$ws.Range("H113").Value = 7730.5757
$ws.Range("I113").Value = 6675.077
$ws.Range("J113").Value = 11651
$ws.Range("K113").Value = 6675.077
$ws.Range("L113").Value = 11651
$ws.Range("M113").Value = -4505.077
$ws.Range("N113").Value = -15991
$ws.Range("H122").Value = 3856.8572
$ws.Range("I122").Value = 2999.8
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 8999.400000000001
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -6549.400000000001
$ws.Range("N122").Value = -22898.5
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040
$ws.Range("H132").Value = 6494825
$ws.Range("I132").Value = 16695141
$ws.Range("K132").Value = 50085423
$ws.Range("M132").Value = -50082893
$ws.Range("H136").Value = 7145592.5
$ws.Range("I136").Value = 7815254
$ws.Range("K136").Value = 23445762
$ws.Range("M136").Value = -23443212

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3474.75
$ws.Range("J96").Value = 3474.75
$ws.Range("L96").Value = 3474.75
$ws.Range("N96").Value = -6220.75
$ws.Range("H107").Value = 5923.4116
$ws.Range("I107").Value = 5763.909
$ws.Range("J107").Value = 6215.8335
$ws.Range("K107").Value = 17291.727
$ws.Range("L107").Value = 18647.5005
$ws.Range("M107").Value = -15371.727
$ws.Range("N107").Value = -22487.5005
$ws.Range("H122").Value = 57995.5
$ws.Range("I122").Value = 1622.7693
$ws.Range("J122").Value = 162687.72
$ws.Range("K122").Value = 4868.3079
$ws.Range("L122").Value = 488063.16
$ws.Range("M122").Value = -2418.3079
$ws.Range("N122").Value = -492963.16
$ws.Range("H125").Value = 55713.285
$ws.Range("I125").Value = 29999
$ws.Range("J125").Value = 59999
$ws.Range("K125").Value = 29999
$ws.Range("L125").Value = 59999
$ws.Range("M125").Value = -25079
$ws.Range("N125").Value = -69839
$ws.Range("H126").Value = 3345.5293
$ws.Range("I126").Value = 3552
$ws.Range("J126").Value = 2674.5
$ws.Range("K126").Value = 10656
$ws.Range("L126").Value = 8023.5
$ws.Range("M126").Value = -8186
$ws.Range("N126").Value = -12963.5
$ws.Range("H131").Value = 46152.75
$ws.Range("J131").Value = 46152.75
$ws.Range("L131").Value = 46152.75
$ws.Range("N131").Value = -56232.75
$ws.Range("H132").Value = 6175849
$ws.Range("I132").Value = 7939499
$ws.Range("J132").Value = 3074.8333
$ws.Range("K132").Value = 23818497
$ws.Range("L132").Value = 9224.499899999999
$ws.Range("M132").Value = -23815967
$ws.Range("N132").Value = -14284.4999
$ws.Range("H136").Value = 10145821
$ws.Range("I136").Value = 2071323.6
$ws.Range("J136").Value = 66667300
$ws.Range("K136").Value = 6213970.800000001
$ws.Range("L136").Value = 200001900
$ws.Range("M136").Value = -6211420.800000001
$ws.Range("N136").Value = -200007000
